# off_white_compare.xlsx update
# - Re-sorts / refreshes the stock-by-store table (rows 2-23) so that it is
#   ordered alphabetically by Title (matches the "aggiunto readme a teerayuth"
#   commit's refreshed export) and brings the "Quantita magazzino" (F) /
#   "Q.ty_Looker" (E) columns up to date.
# - Sets an explicit cell selection (F19) and per-column widths to match the
#   refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.1666666666667
$ws.Columns.Item(2).ColumnWidth = 12.6666666666667
$ws.Columns.Item(3).ColumnWidth = 14.1666666666667
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 12.6666666666667
$ws.Columns.Item(6).ColumnWidth = 13.6666666666667

# --- Data rows (row 1 is the header and is unchanged) -------------------
# Each entry: row, Title, Barcode, Filiale, Variant SKU, Q.ty_Looker, Quantita magazzino
$rows = @(
    @(2,  "Off White Artemisia White",          8052865854535, "0M", "ARTEMISIA 10207", 1, 1),
    @(3,  "Off White Arthur Havana",             8052865804370, "0M", "ARTHUR 26455",    1, 1),
    @(4,  "Off White Arthur Havana",             8052865804370, "D",  "ARTHUR 26455",    1, 1),
    @(5,  "Off White Arthur Black",              8052865781299, "0M", "ARTHUR 11007",    1, 1),
    @(6,  "Off White Baltimore Silver Dark",     8051594080369, "0M", "BALTIMORE 17207", 1, 1),
    @(7,  "Off White Catalina Fuchsia",          8051594080055, "0M", "CATALINA 16807",  1, 1),
    @(8,  "Off White Cady Black",                8052865781251, "D",  "CADY 11007",      2, 2),
    @(9,  "Off White Francisco Black",           8052865854290, "0M", "FRANCISCO 11007", 1, 1),
    @(10, "Off White Catalina Black",            8052865781268, "0M", "CATALINA 11007",  2, 2),
    @(11, "Off White Francisco Blue",            8052865854313, "0M", "FRANCISCO 14507", 1, 1),
    @(12, "Off White Leonardo White",            8051594080062, "0M", "LEONARDO 10107",  1, 1),
    @(13, "Off White Leonardo Black",            8051594257587, "D",  "LEONARDO 11007",  1, 1),
    @(14, "Off White Leonardo Havana",           8051594080079, "D",  "LEONARDO 16055",  1, 1),
    @(15, "Off White Lucio Havana",               8052865854405, "0M", "LUCIO 16055",    1, 1),
    @(16, "Off White Savannah Havana",           8051594079806, "0M", "SAVANNAH 16055",  1, 1),
    @(17, "Off White Seattle Multicolor",        8051594080239, "0M", "SEATTLE 11764",   1, 1),
    @(18, "Off White Seattle Black Multicolor",  8051594080215, "0M", "SEATTLE 11007",   1, 1),
    @(19, "Off White Seattle Black Multicolor",  8051594080215, "D",  "SEATTLE 11007",   1, 1),
    @(20, "Off White Virgil Black/Blue",         8052865804356, "0M", "VIRGIL 21045",    1, 1),
    @(21, "Off White Virgil Black",              8052865781282, "0M", "VIRGIL 11007",    2, 2),
    @(22, "Off White Zurich Black",              8051594030036, "0M", "ZURICH 11007",    1, 1),
    @(23, "Off White Volcanite Black",           8051594376875, "0M", "VOLCANITE 11007", 1, 1)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# --- Selection -----------------------------------------------------------
$ws.Range("F19").Select()
